# Sprint 2 Burndown Chart - update Sprint Backlog / Burndown Chart to its
# latest values. The only underlying data change in this revision is the
# "Actual" figure recorded for Day 6 (row 8): it drops from 9 to 4 on the
# Folha1 worksheet (column C = "Actual"). The burndown chart's "Actual"
# series reads its values straight from Folha1!C3:C9, so updating the cell
# is what drives the chart update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout: A = Day, B = Planned, C = Actual (row 2 holds the headers,
# data starts at row 3). Day 6 is row 8.
$ws.Range("C8").Value = 4
